$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Price" values in column D are stored as text in the source data (e.g. "604.38"),
# even when they look like plain decimals. Pre-mark the cells that would otherwise be
# auto-parsed as numbers with a Text number format so the assigned string is preserved as text.
$numericLookingPriceCells = $ws.Range("D5,D6,D11,D12,D14,D19,D20,D21,D22,D23,D25,D26,D27,D29,D31,D34,D39,D40,D41,D43,D45,D46,D48,D49,D51")
foreach ($area in $numericLookingPriceCells.Areas) {
    $area.NumberFormat = "@"
}

$ws.Range("D2").Value = "65.818.39"
$ws.Range("E2").Value = "  +0.85%  "
$ws.Range("D3").Value = "3.605.60"
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "604.38"
$ws.Range("E5").Value = "  +0.17%  "
$ws.Range("D6").Value = "136.88"
$ws.Range("E6").Value = "  -1.96%  "
$ws.Range("D7").Value = "3.604.01"
$ws.Range("E7").Value = "  +1.79%  "
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("E9").Value = "  +1.83%  "
$ws.Range("E10").Value = "  +0.12%  "
$ws.Range("D11").Value = "7.24"
$ws.Range("E11").Value = "  +4.20%  "
$ws.Range("D12").Value = "0.392"
$ws.Range("E12").Value = "  +0.16%  "
$ws.Range("D13").Value = "4.219.48"
$ws.Range("E13").Value = "  +1.88%  "
$ws.Range("D14").Value = "28.19"
$ws.Range("E14").Value = "  +3.19%  "
$ws.Range("E15").Value = "  +0.43%  "
$ws.Range("D16").Value = "3.606.01"
$ws.Range("E16").Value = "  +1.73%  "
$ws.Range("E17").Value = "  -0.26%  "
$ws.Range("D18").Value = "65.924.32"
$ws.Range("E18").Value = "  +0.88%  "
$ws.Range("D19").Value = "10.09"
$ws.Range("E19").Value = "  -2.58%  "
$ws.Range("D20").Value = "14.74"
$ws.Range("E20").Value = "  +2.52%  "
$ws.Range("D21").Value = "5.91"
$ws.Range("E21").Value = "  -0.73%  "
$ws.Range("D22").Value = "397.97"
$ws.Range("E22").Value = "  +0.58%  "
$ws.Range("D23").Value = "0.592"
$ws.Range("E23").Value = "  +2.90%  "
$ws.Range("D24").Value = "3.752.80"
$ws.Range("E24").Value = "  +1.94%  "
$ws.Range("D25").Value = "74.62"
$ws.Range("E25").Value = "  +0.87%  "
$ws.Range("D26").Value = "0.999"
$ws.Range("E26").Value = "  +0.02%  "
$ws.Range("D27").Value = "0.0000119"
$ws.Range("E27").Value = "  +2.31%  "
$ws.Range("E28").Value = "  +4.71%  "
$ws.Range("D29").Value = "1.67"
$ws.Range("E29").Value = "  +28.60%  "
$ws.Range("E30").Value = "  +4.66%  "
$ws.Range("D31").Value = "8.63"
$ws.Range("E31").Value = "  +3.87%  "
$ws.Range("E32").Value = "  -0.19%  "
$ws.Range("D33").Value = "3.607.82"
$ws.Range("E33").Value = "  +1.53%  "
$ws.Range("D34").Value = "24.59"
$ws.Range("E34").Value = "  +3.23%  "
$ws.Range("E35").Value = "  +1.74%  "
$ws.Range("E37").Value = "  +8.32%  "
$ws.Range("E38").Value = "  +3.43%  "
$ws.Range("D39").Value = "7.09"
$ws.Range("E39").Value = "  +1.38%  "
$ws.Range("D40").Value = "171.51"
$ws.Range("E40").Value = "  +1.50%  "
$ws.Range("D41").Value = "0.0839"
$ws.Range("E41").Value = "  +2.96%  "
$ws.Range("E42").Value = "  +1.79%  "
$ws.Range("D43").Value = "26.24"
$ws.Range("E43").Value = "  -1.76%  "
$ws.Range("E44").Value = "  +0.94%  "
$ws.Range("D45").Value = "1.26"
$ws.Range("E45").Value = "  +4.44%  "
$ws.Range("D46").Value = "4.54"
$ws.Range("E46").Value = "  +1.91%  "
$ws.Range("E47").Value = "  +0.00%  "
$ws.Range("D48").Value = "1.71"
$ws.Range("E48").Value = "  +1.45%  "
$ws.Range("D49").Value = "7.09"
$ws.Range("E49").Value = "  +3.82%  "
$ws.Range("D50").Value = "2.433.86"
$ws.Range("E50").Value = "  -0.88%  "
$ws.Range("D51").Value = "0.0274"
$ws.Range("E51").Value = "  +3.90%  "
